$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Daniel / Ojeda) -- keep name, update sport / department / cars
$ws.Range("C2").Value = "baloncesto"
$ws.Range("E2").Value = "Petén"
$ws.Range("G2").Value = "ford, chrysler, toyota, nissan"

# Row 3 -- replaced placeholder row with Jose Perez
$ws.Range("A3").Value = "Jose "
$ws.Range("B3").Value = "Perez"
$ws.Range("C3").Value = "béisbol"
$ws.Range("D3").Value = "No estoy seguro"
$ws.Range("E3").Value = "Huehuetenango"
$ws.Range("F3").Value = "No"
$ws.Range("G3").Value = "ford, chrysler"

# Row 4 -- replaced placeholder row with Maria Gomez
$ws.Range("A4").Value = "Maria "
$ws.Range("B4").Value = "Gomez"
$ws.Range("C4").Value = "béisbol"
$ws.Range("D4").Value = "Femenino"
$ws.Range("E4").Value = "Petén"
$ws.Range("F4").Value = "Sí"
$ws.Range("G4").Value = "ford, toyota"

# Row 5 -- new row, Jarol Lemus
$ws.Range("A5").Value = "Jarol "
$ws.Range("B5").Value = "Lemus"
$ws.Range("C5").Value = "tenis"
$ws.Range("D5").Value = "Masculino"
$ws.Range("E5").Value = "Chiquimula"
$ws.Range("F5").Value = "Sí"
$ws.Range("G5").Value = "toyota"

# Row 6 -- new row, Ana Ramos
$ws.Range("A6").Value = "Ana "
$ws.Range("B6").Value = "Ramos "
$ws.Range("C6").Value = "béisbol"
$ws.Range("D6").Value = "Femenino"
$ws.Range("E6").Value = "Santa Rosa"
$ws.Range("F6").Value = "Sí"
$ws.Range("G6").Value = "ford"

# Column E got a touch wider to fit "Huehuetenango"
# (ColumnWidth in the COM model is offset ~5/6 from the stored OOXML width,
#  so subtract that offset to land on a stored width of exactly 15)
$ws.Columns("E").ColumnWidth = 15 - (5/6)
